$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 855.2273
$ws.Range("I19").Value = 987.125
$ws.Range("J19").Value = 779.8570999999999
$ws.Range("K19").Value = 987.125
$ws.Range("L19").Value = 779.8570999999999
$ws.Range("M19").Value = -812.125
$ws.Range("N19").Value = -1129.8571

$ws.Range("H40").Value = 888.61536
$ws.Range("I40").Value = 918.75
$ws.Range("J40").Value = 840.4
$ws.Range("K40").Value = 918.75
$ws.Range("L40").Value = 840.4
$ws.Range("M40").Value = -743.75
$ws.Range("N40").Value = -1190.4

$ws.Range("H41").Value = 630.4286
$ws.Range("I41").Value = 699.26666
$ws.Range("K41").Value = 699.26666
$ws.Range("M41").Value = -259.26666

$ws.Range("H43").Value = 74553.734
$ws.Range("I43").Value = 33506.668
$ws.Range("J43").Value = 84815.5
$ws.Range("K43").Value = 33506.668
$ws.Range("L43").Value = 84815.5
$ws.Range("M43").Value = -33437.668
$ws.Range("N43").Value = -84953.5

$ws.Range("H53").Value = 184.16667
$ws.Range("J53").Value = 347.33334
$ws.Range("L53").Value = 347.33334
$ws.Range("N53").Value = -1621.33334

$ws.Range("H70").Value = 3186.2
$ws.Range("J70").Value = 3200.4
$ws.Range("L70").Value = 9601.200000000001
$ws.Range("N70").Value = -10141.2

$ws.Range("H73").Value = 3186.2
$ws.Range("J73").Value = 3200.4
$ws.Range("L73").Value = 9601.200000000001
$ws.Range("N73").Value = -11473.2

$ws.Range("H76").Value = 3460.9375
$ws.Range("I76").Value = 3025.3865
$ws.Range("K76").Value = 3025.3865
$ws.Range("M76").Value = -2710.3865

$ws.Range("H79").Value = 3460.9375
$ws.Range("I79").Value = 3025.3865
$ws.Range("K79").Value = 3025.3865
$ws.Range("M79").Value = -1933.3865

$ws.Range("H141").Value = 2565.1724
$ws.Range("I141").Value = 2088.6843
$ws.Range("J141").Value = 3470.5
$ws.Range("K141").Value = 6266.0529
$ws.Range("L141").Value = 10411.5
$ws.Range("M141").Value = -1086.0529
$ws.Range("N141").Value = -20771.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3058.3333
$ws.Range("I63").Value = 2670
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2670
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1984
$ws.Range("N63").Value = -6372

$ws.Range("H66").Value = 3058.3333
$ws.Range("I66").Value = 2670
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 13350
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -9918
$ws.Range("N66").Value = -31864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 945.2381
$ws.Range("I22").Value = 3616.6667
$ws.Range("K22").Value = 3616.6667
$ws.Range("M22").Value = -3443.6667

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H82").Value = 19713.75
$ws.Range("I82").Value = 9428.5
$ws.Range("J82").Value = 29999
$ws.Range("K82").Value = 9428.5
$ws.Range("L82").Value = 29999
$ws.Range("M82").Value = -9045.5
$ws.Range("N82").Value = -30765

$ws.Range("H85").Value = 19713.75
$ws.Range("I85").Value = 9428.5
$ws.Range("J85").Value = 29999
$ws.Range("K85").Value = 9428.5
$ws.Range("L85").Value = 29999
$ws.Range("M85").Value = -8102.5
$ws.Range("N85").Value = -32651

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H88").Value = 50000
$ws.Range("J88").Value = 50000
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50812

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H91").Value = 50000
$ws.Range("J91").Value = 50000
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -52808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2719.4849
$ws.Range("I58").Value = 962.619
$ws.Range("K58").Value = 962.619
$ws.Range("M58").Value = -759.619

$ws.Range("H74").Value = 11772.25
$ws.Range("J74").Value = 11772.25
$ws.Range("L74").Value = 11772.25
$ws.Range("N74").Value = -13520.25

$ws.Range("H77").Value = 11772.25
$ws.Range("J77").Value = 11772.25
$ws.Range("L77").Value = 35316.75
$ws.Range("N77").Value = -44052.75

$ws.Range("H136").Value = 2719.4849
$ws.Range("I136").Value = 962.619
$ws.Range("K136").Value = 2887.857
$ws.Range("M136").Value = -337.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 149.3
$ws.Range("I12").Value = 232.71428
$ws.Range("J12").Value = 104.38461
$ws.Range("K12").Value = 698.14284
$ws.Range("L12").Value = 313.15383
$ws.Range("M12").Value = -525.14284
$ws.Range("N12").Value = -659.15383

$ws.Range("H33").Value = 1940.8572
$ws.Range("I33").Value = 558
$ws.Range("J33").Value = 7818
$ws.Range("K33").Value = 3348
$ws.Range("L33").Value = 46908
$ws.Range("M33").Value = -3065
$ws.Range("N33").Value = -47474

$ws.Range("H38").Value = 99.166664
$ws.Range("I38").Value = 87.77778000000001
$ws.Range("K38").Value = 263.33334
$ws.Range("M38").Value = 83.66665999999998

$ws.Range("H122").Value = 346.9
$ws.Range("I122").Value = 296.55554
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 2668.99986
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -218.9998599999999
$ws.Range("N122").Value = -12100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2533.2
$ws.Range("I80").Value = 2040
$ws.Range("J80").Value = 2656.5
$ws.Range("K80").Value = 2040
$ws.Range("L80").Value = 2656.5
$ws.Range("M80").Value = -1042
$ws.Range("N80").Value = -4652.5

$ws.Range("H83").Value = 2533.2
$ws.Range("I83").Value = 2040
$ws.Range("J83").Value = 2656.5
$ws.Range("K83").Value = 10200
$ws.Range("L83").Value = 13282.5
$ws.Range("M83").Value = -5208
$ws.Range("N83").Value = -23266.5

$ws.Range("H92").Value = 8050.3335
$ws.Range("J92").Value = 8050.3335
$ws.Range("L92").Value = 8050.3335
$ws.Range("N92").Value = -11794.3335

$ws.Range("H113").Value = 1582.7567
$ws.Range("J113").Value = 1388.6666
$ws.Range("L113").Value = 1388.6666
$ws.Range("N113").Value = -5728.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 800
$ws.Range("I16").Value = 625
$ws.Range("K16").Value = 625
$ws.Range("M16").Value = -455

$ws.Range("H22").Value = 1250.2222
$ws.Range("I22").Value = 725
$ws.Range("J22").Value = 1670.4
$ws.Range("K22").Value = 725
$ws.Range("L22").Value = 1670.4
$ws.Range("M22").Value = -430
$ws.Range("N22").Value = -2260.4

$ws.Range("H27").Value = 1250.2222
$ws.Range("I27").Value = 725
$ws.Range("J27").Value = 1670.4
$ws.Range("K27").Value = 725
$ws.Range("L27").Value = 1670.4
$ws.Range("M27").Value = -618
$ws.Range("N27").Value = -1884.4

$ws.Range("H64").Value = 18362.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 18362.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 18362.5
$ws.Range("N64").Value = -18812.5
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 18362.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 18362.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 18362.5
$ws.Range("N67").Value = -19922.5
$ws.Range("M67").ClearContents()

$ws.Range("H93").Value = 1291.1
$ws.Range("I93").Value = 948.7143
$ws.Range("J93").Value = 2090
$ws.Range("K93").Value = 948.7143
$ws.Range("L93").Value = 2090
$ws.Range("M93").Value = 299.2857
$ws.Range("N93").Value = -4586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 15666.667
$ws.Range("J63").Value = 15666.667
$ws.Range("L63").Value = 15666.667
$ws.Range("N63").Value = -16914.667

$ws.Range("H66").Value = 15666.667
$ws.Range("J66").Value = 15666.667
$ws.Range("L66").Value = 47000.001
$ws.Range("N66").Value = -53240.001

$ws.Range("H80").Value = 42285.715
$ws.Range("J80").Value = 42285.715
$ws.Range("L80").Value = 42285.715
$ws.Range("N80").Value = -44281.715

$ws.Range("H81").Value = 1432.8928
$ws.Range("I81").Value = 897.44446
$ws.Range("J81").Value = 1686.5264
$ws.Range("K81").Value = 1794.88892
$ws.Range("L81").Value = 3373.0528
$ws.Range("M81").Value = -733.8889200000001
$ws.Range("N81").Value = -5495.052799999999

$ws.Range("H83").Value = 42285.715
$ws.Range("J83").Value = 42285.715
$ws.Range("L83").Value = 126857.145
$ws.Range("N83").Value = -136841.145

$ws.Range("H84").Value = 1432.8928
$ws.Range("I84").Value = 897.44446
$ws.Range("J84").Value = 1686.5264
$ws.Range("K84").Value = 8974.444600000001
$ws.Range("L84").Value = 16865.264
$ws.Range("M84").Value = -3670.444600000001
$ws.Range("N84").Value = -27473.264
